$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows starting at row 25 (old rows 25-31 shift down to 33-39).
$ws.Rows("25:32").Insert()

# Row 24 - "Rotate Image" entry renamed to "RotateImage" (date/platform unchanged)
$ws.Range("C24").Value = "RotateImage"

# Row 25 (new) - RotateImage2 / IntelliJ IDEA (PC)
$ws.Range("C25").Value = "RotateImage2"
$ws.Range("E25").Value = "IntelliJ IDEA (PC)"

# Row 26 (new) - Rotate Image (Solution 1) / LeetCode
$ws.Range("C26").Value = "Rotate Image (Solution 1)"
$ws.Range("E26").Value = "LeetCode"

# Row 27 (new) - Rotate Image (Solution 2) / LeetCode
$ws.Range("C27").Value = "Rotate Image (Solution 2)"
$ws.Range("E27").Value = "LeetCode"

# Row 28 (new) - Flipping an Image / Bosscoder Academy
$ws.Range("C28").Value = "Flipping an Image"
$ws.Range("E28").Value = "Bosscoder Academy"

# Row 29 (new) - Set Matrix Zero (Solution 1) / LeetCode
$ws.Range("C29").Value = "Set Matrix Zero (Solution 1)"
$ws.Range("E29").Value = "LeetCode"

# Row 30 (new) - Date 9, Set Matrix Zero (Solution 2) / LeetCode
$ws.Range("B30").Value = 9
$ws.Range("C30").Value = "Set Matrix Zero (Solution 2)"
$ws.Range("E30").Value = "LeetCode"

# Row 31 (new) - Reshape the Matrix / LeetCode
$ws.Range("C31").Value = "Reshape the Matrix"
$ws.Range("E31").Value = "LeetCode"

# Row 32 (new) - Date 10, FillMatrixByRecursion / IntelliJ IDEA (PC)
$ws.Range("B32").Value = 10
$ws.Range("C32").Value = "FillMatrixByRecursion"
$ws.Range("E32").Value = "IntelliJ IDEA (PC)"

# Row 33 (was old row 25) - now PermutationsOfAString, platform unchanged (IntelliJ)
$ws.Range("C33").Value = "PermutationsOfAString"

# Row 34 (was old row 26) - Date 11, SubsetsOfAString / IntelliJ IDEA (PC)
$ws.Range("B34").Value = 11
$ws.Range("C34").Value = "SubsetsOfAString"
$ws.Range("E34").Value = "IntelliJ IDEA (PC)"

# Row 35 (was old row 27) - Subsets (Solution 1) / LeetCode
$ws.Range("C35").Value = "Subsets (Solution 1)"
$ws.Range("E35").Value = "LeetCode"

# Row 36 (was old row 28) - Subsets (Solution 2) / LeetCode
$ws.Range("C36").Value = "Subsets (Solution 2)"

# Row 37 (was old row 29) - Date 12, Permutations / Bosscoder Academy
$ws.Range("B37").Value = 12
$ws.Range("C37").Value = "Permutations"
$ws.Range("E37").Value = "Bosscoder Academy"

# Row 38 (was old row 30) - clear leftover Date, Permutations / LeetCode
$ws.Range("B38").ClearContents()
$ws.Range("C38").Value = "Permutations"
$ws.Range("E38").Value = "LeetCode"

# Row 39 (was old row 31) - Date 13, clear leftover Problem/Platform
$ws.Range("B39").Value = 13
$ws.Range("C39").ClearContents()
$ws.Range("E39").ClearContents()

# Update the sheet view to match the author's final cursor position
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C39").Select()
